# Updated symbol list on Wed Dec 21 09:58:10 UTC 2022 with GitHub Actions
#
# Applies the latest coinranking.com price refresh to the "cryptos" sheet:
#   - refreshed Price (column D) values for the rows whose quote moved
#   - rows 41-43 reshuffled (BKEXToken, CEJI, KickToken rotated into new
#     rank order) with refreshed Coin / Link / Price / Volume(1h) cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing genuine Text storage (these price
# cells are stored as text, e.g. t="inlineStr"/shared-string - not
# numbers), then drop the temporary "Text" number-format override so the
# cell's style index is left exactly as it was before the edit.
function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---- Column D (Price) refreshes ----------------------------------------
$priceChanges = @{
    "D2"  = "248.89"
    "D3"  = "22.58"
    "D4"  = "5.395"
    "D5"  = "0.05692"
    "D6"  = "3.404"
    "D7"  = "6.325"
    "D8"  = "0.8054"
    "D9"  = "0.9162"
    "D10" = "0.1401"
    "D11" = "0.07438"
    "D12" = "0.03134"
    "D13" = "0.03033"
    "D14" = "0.09378"
    "D16" = "0.001573"
    "D17" = "0.04785"
    "D18" = "0.01826"
    "D19" = "0.0005845"
    "D21" = "0.004997"
    "D22" = "0.001005"
    "D24" = "3.702"
    "D25" = "2.200"
    "D40" = "0.04018"
    "D44" = "0.007570"
    "D45" = "0.00005787"
    "D46" = "0.00000000749"
    "D47" = "0.4986"
    "D48" = "0.2101"
    "D49" = "0.00002098"
    "D50" = "0.01009"
}

foreach ($addr in $priceChanges.Keys) {
    Set-TextValue $ws $addr $priceChanges[$addr]
}

# ---- Rows 41-43 rank reshuffle ------------------------------------------
# BKEXToken / CEJI / KickToken rotate into new rank positions with
# refreshed Coin, Link, Price and Volume(1h) values.
$rowUpdates = @{
    41 = @("BKEXToken", "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk", "0.1071", "40BKEXTokenBKK")
    42 = @("CEJI", "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji", "0.002719", "41CEJICEJI")
    43 = @("KickToken", "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick", "0.006849", "42KickTokenKICK")
}

foreach ($row in $rowUpdates.Keys) {
    $vals = $rowUpdates[$row]
    $addrB = "B" + $row
    $addrC = "C" + $row
    $addrD = "D" + $row
    $addrE = "E" + $row
    Set-TextValue $ws $addrB $vals[0]
    Set-TextValue $ws $addrC $vals[1]
    Set-TextValue $ws $addrD $vals[2]
    Set-TextValue $ws $addrE $vals[3]
}
